# Fix spreadsheets, clocks.c, and trace/tools/ scripts for 250Mhz operation.
# The workbook's "IOPLL" divisor (C28, C30 on Sheet1) changes from 8 to 6,
# which cascades through the dependent formulas across the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C28").Value = 6
$ws.Range("C30").Value = 6

# Update the active selection on Sheet1 to match the saved view state.
$ws.Activate() | Out-Null
$ws.Range("I4").Select() | Out-Null
